# Auto-generated edit script for lipidcane_spearman_-2_agile.xlsx
# Adds new parameter rows (15-23) to the Spearman sensitivity sheet and
# refreshes the sensitivity coefficients for the existing rows (4-14).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Refresh sensitivity coefficients for existing parameters (rows 4-14) ---
$ws.Range("C4").Value = -0.01872018583680743
$ws.Range("E4").Value = 0.008664924826596991
$ws.Range("G4").Value = 0.003205590944223637
$ws.Range("H4").Value = -0.02225724175428967
$ws.Range("I4").Value = -0.01926712281868491
$ws.Range("J4").Value = 0.01462343566151529
$ws.Range("C5").Value = -0.008663906458556258
$ws.Range("E5").Value = -0.01244344091373763
$ws.Range("G5").Value = -0.009917670636706823
$ws.Range("H5").Value = 0.002330263197210527
$ws.Range("I5").Value = -0.0009816489032659559
$ws.Range("J5").Value = 0.00281660288534416
$ws.Range("C6").Value = 0.005172327566893102
$ws.Range("E6").Value = 0.01592268927690757
$ws.Range("G6").Value = 0.01050879555635182
$ws.Range("H6").Value = 0.0002081100563244022
$ws.Range("I6").Value = -0.003178638271145531
$ws.Range("J6").Value = 0.01306171845696569
$ws.Range("C7").Value = 0.003711844660473786
$ws.Range("E7").Value = 0.00799635775985431
$ws.Range("G7").Value = 0.01061265152850606
$ws.Range("H7").Value = 0.009266903602676142
$ws.Range("I7").Value = 0.0009264983410599334
$ws.Range("J7").Value = -0.008311832983149059
$ws.Range("C8").Value = 0.09906571183462846
$ws.Range("E8").Value = 0.003319033092761323
$ws.Range("G8").Value = -0.01287873382714935
$ws.Range("H8").Value = 0.9813190557007622
$ws.Range("I8").Value = 0.974378954879158
$ws.Range("J8").Value = 0.01614513631948846
$ws.Range("C9").Value = 0.9542819602192784
$ws.Range("E9").Value = 0.004938495365539814
$ws.Range("G9").Value = -0.001715374340614973
$ws.Range("H9").Value = -0.02666081367443255
$ws.Range("I9").Value = -0.02604824427392977
$ws.Range("J9").Value = 0.02174168477581664
$ws.Range("C10").Value = 0.004633077977323118
$ws.Range("E10").Value = -0.004318801804752071
$ws.Range("G10").Value = 0.0002164983446599338
$ws.Range("H10").Value = 0.01139245264769811
$ws.Range("I10").Value = 0.008461169330446773
$ws.Range("J10").Value = -0.01603021816808957
$ws.Range("C11").Value = -0.0139813374552535
$ws.Range("E11").Value = 0.007653537618141504
$ws.Range("G11").Value = 0.00455106978204279
$ws.Range("H11").Value = 0.007043801081752042
$ws.Range("I11").Value = 0.00758259678330387
$ws.Range("J11").Value = 0.001858249079641509
$ws.Range("C12").Value = 0.02768172667526906
$ws.Range("E12").Value = 0.005015854472634178
$ws.Range("G12").Value = 0.003109101052364041
$ws.Range("H12").Value = 0.0219217525568701
$ws.Range("I12").Value = 0.0207861575674463
$ws.Range("J12").Value = 0.00850402782752333
$ws.Range("C13").Value = 0.07419288200771526
$ws.Range("E13").Value = 0.004447120785884831
$ws.Range("G13").Value = -0.00359526350381054
$ws.Range("H13").Value = -0.01190910518036421
$ws.Range("I13").Value = 0.2067890159195606
$ws.Range("J13").Value = 0.007917316829748527
$ws.Range("C14").Value = -0.181145151405806
$ws.Range("E14").Value = -0.001062775434511017
$ws.Range("G14").Value = -0.001173795694951828
$ws.Range("H14").Value = 0.01941646359265854
$ws.Range("I14").Value = 0.02511686980467479
$ws.Range("J14").Value = 0.001909759466957954

# --- New parameter rows (15-23): write column A (Element) first, across all new rows ---
$ws.Range("A15").Value = "Stream-crude glycerol"
$ws.Range("A16").Value = "Stream-pure glycerine"
$ws.Range("A17").Value = "Stream-cellulase"
$ws.Range("A19").Value = "Pretreatment reactor system"
$ws.Range("A20").Value = "Pretreatment and saccharification"
$ws.Range("A22").Value = "Cofermenation"

# --- New parameter rows (15-23): write column B (Parameter) next, across all new rows ---
$ws.Range("B15").Value = "Price [USD/kg]"
$ws.Range("B16").Value = "Price [USD/kg]"
$ws.Range("B17").Value = "Price [USD/kg]"
$ws.Range("B18").Value = "Cellulase loading [wt. % cellulose]"
$ws.Range("B19").Value = "Base cost [million USD]"
$ws.Range("B20").Value = "Glucose yield [%]"
$ws.Range("B21").Value = "Xylose yield [%]"
$ws.Range("B22").Value = "Glucose to ethanol yield [%]"
$ws.Range("B23").Value = "Xylose to ethanol yield [%]"

# --- New parameter rows (15-23): sensitivity coefficients ---
$ws.Range("C15").Value = -0.01277722947108918
$ws.Range("E15").Value = 0.009341105845644232
$ws.Range("G15").Value = 0.01233837332553493
$ws.Range("H15").Value = 0.002626834857073394
$ws.Range("I15").Value = 0.003835622073424882
$ws.Range("J15").Value = -0.001807113956349225
$ws.Range("C16").Value = -0.004111017092440683
$ws.Range("E16").Value = -0.001598739519949581
$ws.Range("G16").Value = -0.0009761016390440654
$ws.Range("H16").Value = 0.02953779814151192
$ws.Range("I16").Value = 0.02279120855964833
$ws.Range("J16").Value = 0.0004208061870931857
$ws.Range("C17").Value = -0.05054158214966328
$ws.Range("E17").Value = -0.004070838786833551
$ws.Range("G17").Value = 0.001633359521334381
$ws.Range("H17").Value = 0.04318200163128006
$ws.Range("I17").Value = 0.04670276855611073
$ws.Range("J17").Value = -0.01451312349441295
$ws.Range("C18").Value = 0.02603518644940745
$ws.Range("E18").Value = -0.01147393639495745
$ws.Range("G18").Value = -0.006629770729190828
$ws.Range("H18").Value = 0.01299372944774918
$ws.Range("I18").Value = 0.01622840628113625
$ws.Range("J18").Value = -0.003996017681279603
$ws.Range("C19").Value = -0.03052129302885172
$ws.Range("E19").Value = 0.02988183460327338
$ws.Range("G19").Value = 0.03090082405203296
$ws.Range("H19").Value = 0.1878021709680868
$ws.Range("I19").Value = 0.003863478394539135
$ws.Range("J19").Value = -0.00146147241409544
$ws.Range("C20").Value = 0.03479382695975308
$ws.Range("E20").Value = 0.2472170067846803
$ws.Range("G20").Value = 0.232911876388475
$ws.Range("H20").Value = 0.01166759240270369
$ws.Range("I20").Value = -0.001004246920169877
$ws.Range("J20").Value = -0.1463784983204333
$ws.Range("C21").Value = 0.08487786723511467
$ws.Range("E21").Value = 0.3982084093363363
$ws.Range("G21").Value = 0.09530755370030213
$ws.Range("H21").Value = 0.006564748294589931
$ws.Range("I21").Value = 0.02428224865128994
$ws.Range("J21").Value = 0.9343579361897016
$ws.Range("C22").Value = 0.01001100635244025
$ws.Range("E22").Value = -0.008940768837630753
$ws.Range("G22").Value = -0.006917342292693691
$ws.Range("H22").Value = -0.005763718118548724
$ws.Range("I22").Value = -0.0009532179261287168
$ws.Range("J22").Value = -0.006714976224337031
$ws.Range("C23").Value = 0.0822693774027751
$ws.Range("E23").Value = 0.8707920958876838
$ws.Range("G23").Value = 0.9599115147804604
$ws.Range("H23").Value = 0.01177128565485143
$ws.Range("I23").Value = -0.01008051323522053
$ws.Range("J23").Value = -0.2215270908548003

# --- Copy cell formatting (bold, border, centered/top aligned) onto the new cells ---
$ws.Range("A14:B14").Copy()
$ws.Range("A15:B23").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Merge the Element column for grouped parameter rows ---
$ws.Range("A17:A18").Merge()
$ws.Range("A20:A21").Merge()
$ws.Range("A22:A23").Merge()

